$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string: "clinical trial search app" (index 27, 0-based)
# New rows appended to the worklog: rows 55 and 56

$ws.Range("A54").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("A56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A55").Value = 43132
$ws.Range("B55").Value = "clinical trial search app"

$ws.Range("A56").Value = 43133
$ws.Range("B56").Value = "clinical trial search app"

$ws.Range("B56").Select()
